$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (20) comparing to a household name: iShares EM LibertyQ (EM) fund.
$ws.Range("A20").Value = "EM LibertyQ (EM)"
$ws.Range("B20").Value = 709646
$ws.Range("C20").Value = "EM"
$ws.Range("E20").Value = "LibertyQ"
$ws.Range("G20").Value = 0.45

# Move the selection to where the user clicked next, as recorded in the workbook.
$ws.Range("E24").Select()
